$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 2 new rows at row 33: new Hotel Info "no information found" test cases ---
$ws.Rows("33:34").Insert()

# --- Insert 21 new rows after row 40 (the shifted former last row): new "invalid/missing token" test cases ---
$ws.Rows("41:61").Insert()

# Row 33 (Test Case #32)
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = '5.1.0'
$ws.Range("C33").Value = 'The Hotel Info Resource should inform the user if no information was found'
$ws.Range("C33").WrapText = $true
$ws.Range("D33").Value = 'Get request to Hotel Info Resource with invalid ID input'
$ws.Range("D33").WrapText = $true
$ws.Range("E33").Value = 'Message stating no information was found'
$ws.Range("E33").WrapText = $true
$ws.Rows(33).RowHeight = 45

# Row 34 (Test Case #33)
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = '5.1.0'
$ws.Range("C34").Value = 'The Hotel Info Resource should inform the user if no information was found'
$ws.Range("C34").WrapText = $true
$ws.Range("D34").Value = 'Get request to Hotel Info Resource with missing ID input'
$ws.Range("D34").WrapText = $true
$ws.Range("E34").Value = 'Message stating no information was found'
$ws.Range("E34").WrapText = $true
$ws.Rows(34).RowHeight = 45

# Row 35 (Test Case #34)
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = '6.0.0
6.1.0
6.2.0'
$ws.Range("B35").WrapText = $true
$ws.Range("C35").Value = 'The user shall register to use the API by creating a username and password'
$ws.Range("C35").WrapText = $true
$ws.Range("D35").Value = 'Post request to User Registration with inputs of username and password'
$ws.Range("D35").WrapText = $true
$ws.Range("E35").Value = 'Message stating ''username'' was created'
$ws.Range("E35").WrapText = $true
$ws.Rows(35).RowHeight = 60

# Row 36 (Test Case #35)
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = '6.0.0'
$ws.Range("C36").Value = 'The user should be informed if there was an error in registration'
$ws.Range("C36").WrapText = $true
$ws.Range("D36").Value = 'Post request to User Registration with invalid inputs'
$ws.Range("D36").WrapText = $true
$ws.Range("E36").Value = 'Error message'
$ws.Range("E36").WrapText = $true
$ws.Rows(36).RowHeight = 45

# Row 37 (Test Case #36)
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = '6.1.1
6.1.2'
$ws.Range("B37").WrapText = $true
$ws.Range("C37").Value = 'The username shall be unique'
$ws.Range("C37").WrapText = $true
$ws.Range("D37").Value = 'Post request to User Registration with an existing username as input'
$ws.Range("D37").WrapText = $true
$ws.Range("E37").Value = 'Message stating ''username'' already exists'
$ws.Range("E37").WrapText = $true
$ws.Rows(37).RowHeight = 60

# Row 38 (Test Case #37)
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = '7.0.0
8.0.0'
$ws.Range("B38").WrapText = $true
$ws.Range("C38").Value = 'The user shall be required to login to use the API'
$ws.Range("C38").WrapText = $true
$ws.Range("D38").Value = 'Post request to User Login with valid username and password'
$ws.Range("D38").WrapText = $true
$ws.Range("E38").Value = 'Access token created'
$ws.Range("E38").WrapText = $true
$ws.Rows(38).RowHeight = 60

# Row 39 (Test Case #38)
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = '7.1.0'
$ws.Range("B39").WrapText = $true
$ws.Range("C39").Value = 'The user shall be informed if credentials to login are invalid'
$ws.Range("C39").WrapText = $true
$ws.Range("D39").Value = 'Post request to User Login with invalid username and password'
$ws.Range("D39").WrapText = $true
$ws.Range("E39").Value = 'Message stating credentials are invalid'
$ws.Range("E39").WrapText = $true
$ws.Rows(39).RowHeight = 60

# Row 40 (Test Case #39)
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = '8.1.0
8.1.1'
$ws.Range("B40").WrapText = $true
$ws.Range("C40").Value = 'The access token shall expire after 15 minutes'
$ws.Range("C40").WrapText = $true
$ws.Range("D40").Value = 'Get request to Weather Resource with no added input'
$ws.Range("D40").WrapText = $true
$ws.Range("E40").Value = 'Message stating the token expired'
$ws.Range("E40").WrapText = $true
$ws.Rows(40).RowHeight = 45

# Row 41 (Test Case #40)
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = '8.1.0
8.1.1'
$ws.Range("B41").WrapText = $true
$ws.Range("C41").Value = 'The access token shall expire after 15 minutes'
$ws.Range("C41").WrapText = $true
$ws.Range("D41").Value = 'Get request to WeatherFiveDay Resource with no added input'
$ws.Range("D41").WrapText = $true
$ws.Range("E41").Value = 'Message stating the token expired'
$ws.Range("E41").WrapText = $true
$ws.Rows(41).RowHeight = 60

# Row 42 (Test Case #41)
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = '8.1.0
8.1.1'
$ws.Range("B42").WrapText = $true
$ws.Range("C42").Value = 'The access token shall expire after 15 minutes'
$ws.Range("C42").WrapText = $true
$ws.Range("D42").Value = 'Get request to Restaurant Resource with no added input'
$ws.Range("D42").WrapText = $true
$ws.Range("E42").Value = 'Message stating the token expired'
$ws.Range("E42").WrapText = $true
$ws.Rows(42).RowHeight = 45

# Row 43 (Test Case #42)
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = '8.1.0
8.1.1'
$ws.Range("B43").WrapText = $true
$ws.Range("C43").Value = 'The access token shall expire after 15 minutes'
$ws.Range("C43").WrapText = $true
$ws.Range("D43").Value = 'Get request to Event Resource with no added input'
$ws.Range("D43").WrapText = $true
$ws.Range("E43").Value = 'Message stating the token expired'
$ws.Range("E43").WrapText = $true
$ws.Rows(43).RowHeight = 45

# Row 44 (Test Case #43)
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = '8.1.0
8.1.1'
$ws.Range("B44").WrapText = $true
$ws.Range("C44").Value = 'The access token shall expire after 15 minutes'
$ws.Range("C44").WrapText = $true
$ws.Range("D44").Value = 'Get request to Hotel Resource with no added input'
$ws.Range("D44").WrapText = $true
$ws.Range("E44").Value = 'Message stating the token expired'
$ws.Range("E44").WrapText = $true
$ws.Rows(44).RowHeight = 45

# Row 45 (Test Case #44)
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = '8.1.0
8.1.1'
$ws.Range("B45").WrapText = $true
$ws.Range("C45").Value = 'The access token shall expire after 15 minutes'
$ws.Range("C45").WrapText = $true
$ws.Range("D45").Value = 'Get request to Weather Resource with zipcode input'
$ws.Range("D45").WrapText = $true
$ws.Range("E45").Value = 'Message stating the token expired'
$ws.Range("E45").WrapText = $true
$ws.Rows(45).RowHeight = 45

# Row 46 (Test Case #45)
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = '8.1.0
8.1.1'
$ws.Range("B46").WrapText = $true
$ws.Range("C46").Value = 'The access token shall expire after 15 minutes'
$ws.Range("C46").WrapText = $true
$ws.Range("D46").Value = 'Get request to WeatherFiveDay Resource with zipcode input'
$ws.Range("D46").WrapText = $true
$ws.Range("E46").Value = 'Message stating the token expired'
$ws.Range("E46").WrapText = $true
$ws.Rows(46).RowHeight = 60

# Row 47 (Test Case #46)
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = '8.1.0
8.1.1'
$ws.Range("B47").WrapText = $true
$ws.Range("C47").Value = 'The access token shall expire after 15 minutes'
$ws.Range("C47").WrapText = $true
$ws.Range("D47").Value = 'Get request to Restaurant Resource with zipcode input'
$ws.Range("D47").WrapText = $true
$ws.Range("E47").Value = 'Message stating the token expired'
$ws.Range("E47").WrapText = $true
$ws.Rows(47).RowHeight = 45

# Row 48 (Test Case #47)
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = '8.1.0
8.1.1'
$ws.Range("B48").WrapText = $true
$ws.Range("C48").Value = 'The access token shall expire after 15 minutes'
$ws.Range("C48").WrapText = $true
$ws.Range("D48").Value = 'Get request to Event Resource with zipcode input'
$ws.Range("D48").WrapText = $true
$ws.Range("E48").Value = 'Message stating the token expired'
$ws.Range("E48").WrapText = $true
$ws.Rows(48).RowHeight = 45

# Row 49 (Test Case #48)
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = '8.1.0
8.1.1'
$ws.Range("B49").WrapText = $true
$ws.Range("C49").Value = 'The access token shall expire after 15 minutes'
$ws.Range("C49").WrapText = $true
$ws.Range("D49").Value = 'Get request to Hotel Resource with zipcode input'
$ws.Range("D49").WrapText = $true
$ws.Range("E49").Value = 'Message stating the token expired'
$ws.Range("E49").WrapText = $true
$ws.Rows(49).RowHeight = 45

# Row 50 (Test Case #49)
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = '8.1.0
8.1.1'
$ws.Range("B50").WrapText = $true
$ws.Range("C50").Value = 'The access token shall expire after 15 minutes'
$ws.Range("C50").WrapText = $true
$ws.Range("D50").Value = 'Get request to Hotel Info Resource with XID input'
$ws.Range("D50").WrapText = $true
$ws.Range("E50").Value = 'Message stating the token expired'
$ws.Range("E50").WrapText = $true
$ws.Rows(50).RowHeight = 45

# Row 51 (Test Case #50)
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = '8.3.0'
$ws.Range("B51").WrapText = $true
$ws.Range("C51").Value = 'The user shall be informed of missing or invalid access token when making a request to Weather Resource'
$ws.Range("C51").WrapText = $true
$ws.Range("D51").Value = 'Get request to Weather Resource with no added input'
$ws.Range("D51").WrapText = $true
$ws.Range("E51").Value = 'Message indicating invalid access token'
$ws.Range("E51").WrapText = $true
$ws.Rows(51).RowHeight = 60

# Row 52 (Test Case #51)
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = '8.3.0'
$ws.Range("B52").WrapText = $true
$ws.Range("C52").Value = 'The user shall be informed of missing or invalid access token when making a request to WeatherFiveDay Resource'
$ws.Range("C52").WrapText = $true
$ws.Range("D52").Value = 'Get request to WeatherFiveDay Resource with no added input'
$ws.Range("D52").WrapText = $true
$ws.Range("E52").Value = 'Message indicating invalid access token'
$ws.Range("E52").WrapText = $true
$ws.Rows(52).RowHeight = 60

# Row 53 (Test Case #52)
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = '8.3.0'
$ws.Range("B53").WrapText = $true
$ws.Range("C53").Value = 'The user shall be informed of missing or invalid access token when making a request to Restaurant Resource'
$ws.Range("C53").WrapText = $true
$ws.Range("D53").Value = 'Get request to Restaurant Resource with no added input'
$ws.Range("D53").WrapText = $true
$ws.Range("E53").Value = 'Message indicating invalid access token'
$ws.Range("E53").WrapText = $true
$ws.Rows(53).RowHeight = 60

# Row 54 (Test Case #53)
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = '8.3.0'
$ws.Range("B54").WrapText = $true
$ws.Range("C54").Value = 'The user shall be informed of missing or invalid access token when making a request to Event Resource'
$ws.Range("C54").WrapText = $true
$ws.Range("D54").Value = 'Get request to Event Resource with no added input'
$ws.Range("D54").WrapText = $true
$ws.Range("E54").Value = 'Message indicating invalid access token'
$ws.Range("E54").WrapText = $true
$ws.Rows(54).RowHeight = 60

# Row 55 (Test Case #54)
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = '8.3.0'
$ws.Range("B55").WrapText = $true
$ws.Range("C55").Value = 'The user shall be informed of missing or invalid access token when making a request to Hotel Resource'
$ws.Range("C55").WrapText = $true
$ws.Range("D55").Value = 'Get request to Hotel Resource with no added input'
$ws.Range("D55").WrapText = $true
$ws.Range("E55").Value = 'Message indicating invalid access token'
$ws.Range("E55").WrapText = $true
$ws.Rows(55).RowHeight = 60

# Row 56 (Test Case #55)
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = '8.3.0'
$ws.Range("B56").WrapText = $true
$ws.Range("C56").Value = 'The user shall be informed of missing or invalid access token when making a request to Weather Resource'
$ws.Range("C56").WrapText = $true
$ws.Range("D56").Value = 'Get request to Weather Resource with zipcode input'
$ws.Range("D56").WrapText = $true
$ws.Range("E56").Value = 'Message indicating invalid access token'
$ws.Range("E56").WrapText = $true
$ws.Rows(56).RowHeight = 60

# Row 57 (Test Case #56)
$ws.Range("A57").Value = 56
$ws.Range("B57").Value = '8.3.0'
$ws.Range("B57").WrapText = $true
$ws.Range("C57").Value = 'The user shall be informed of missing or invalid access token when making a request to WeatherFiveDay Resource'
$ws.Range("C57").WrapText = $true
$ws.Range("D57").Value = 'Get request to WeatherFiveDay Resource with zipcode input'
$ws.Range("D57").WrapText = $true
$ws.Range("E57").Value = 'Message indicating invalid access token'
$ws.Range("E57").WrapText = $true
$ws.Rows(57).RowHeight = 60

# Row 58 (Test Case #57)
$ws.Range("A58").Value = 57
$ws.Range("B58").Value = '8.3.0'
$ws.Range("B58").WrapText = $true
$ws.Range("C58").Value = 'The user shall be informed of missing or invalid access token when making a request to Restaurant Resource'
$ws.Range("C58").WrapText = $true
$ws.Range("D58").Value = 'Get request to Restaurant Resource with zipcode input'
$ws.Range("D58").WrapText = $true
$ws.Range("E58").Value = 'Message indicating invalid access token'
$ws.Range("E58").WrapText = $true
$ws.Rows(58).RowHeight = 60

# Row 59 (Test Case #58)
$ws.Range("A59").Value = 58
$ws.Range("B59").Value = '8.3.0'
$ws.Range("B59").WrapText = $true
$ws.Range("C59").Value = 'The user shall be informed of missing or invalid access token when making a request to Event Resource'
$ws.Range("C59").WrapText = $true
$ws.Range("D59").Value = 'Get request to Event Resource with zipcode input'
$ws.Range("D59").WrapText = $true
$ws.Range("E59").Value = 'Message indicating invalid access token'
$ws.Range("E59").WrapText = $true
$ws.Rows(59).RowHeight = 60

# Row 60 (Test Case #59)
$ws.Range("A60").Value = 59
$ws.Range("B60").Value = '8.3.0'
$ws.Range("B60").WrapText = $true
$ws.Range("C60").Value = 'The user shall be informed of missing or invalid access token when making a request to Hotel Resource'
$ws.Range("C60").WrapText = $true
$ws.Range("D60").Value = 'Get request to Hotel Resource with zipcode input'
$ws.Range("D60").WrapText = $true
$ws.Range("E60").Value = 'Message indicating invalid access token'
$ws.Range("E60").WrapText = $true
$ws.Rows(60).RowHeight = 60

# Row 61 (Test Case #60)
$ws.Range("A61").Value = 60
$ws.Range("B61").Value = '8.3.0'
$ws.Range("B61").WrapText = $true
$ws.Range("C61").Value = 'The user shall be informed of missing or invalid access token when making a request to Hotel Resource'
$ws.Range("C61").WrapText = $true
$ws.Range("D61").Value = 'Get request to Hotel Info Resource with XID input'
$ws.Range("D61").WrapText = $true
$ws.Range("E61").Value = 'Message indicating invalid access token'
$ws.Range("E61").WrapText = $true
$ws.Rows(61).RowHeight = 60

# --- Update sheet view: scroll position & selection to match the edited state ---
$ws.Range("C61").Select()
$ws.Application.ActiveWindow.ScrollRow = 49

